# Apply updates to both sheets of the workbook: new random A1:A10 values,
# swap the average/std label & formula columns (formula now in column A,
# label now in column B), and update the selected cell on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: misclassification_rates
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("misclassification_rates")

$ws1.Range("A1").Value = 0.184
$ws1.Range("A2").Value = 0.154
$ws1.Range("A3").Value = 0.19400000000000001
$ws1.Range("A4").Value = 0.16200000000000001
$ws1.Range("A5").Value = 0.158
$ws1.Range("A6").Value = 0.182
$ws1.Range("A7").Value = 0.14399999999999999
$ws1.Range("A8").Value = 0.214
$ws1.Range("A9").Value = 0.16
$ws1.Range("A10").Value = 0.17799999999999999

$ws1.Range("A12").Value = $null
$ws1.Range("A13").Value = $null
$ws1.Range("A12").Formula = "=AVERAGE(A1:A10)"
$ws1.Range("B12").Value = "avg"
$ws1.Range("A13").Formula = "=_xlfn.STDEV.S(A1:A10)"
$ws1.Range("B13").Value = "std"

$ws1.Range("A14").Select()

# ---------------------------------------------------------------
# Sheet 2: brier_scores
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("brier_scores")

$ws2.Range("A1").Value = 0.27806748799999997
$ws2.Range("A2").Value = 0.26949211200000001
$ws2.Range("A3").Value = 0.29549150400000002
$ws2.Range("A4").Value = 0.25762707200000001
$ws2.Range("A5").Value = 0.27116889599999999
$ws2.Range("A6").Value = 0.29406406400000001
$ws2.Range("A7").Value = 0.26478438399999998
$ws2.Range("A8").Value = 0.30042103999999997
$ws2.Range("A9").Value = 0.27671766399999997
$ws2.Range("A10").Value = 0.28371555199999998

$ws2.Range("A12").Value = $null
$ws2.Range("A13").Value = $null
$ws2.Range("A12").Formula = "=AVERAGE(A1:A10)"
$ws2.Range("B12").Value = "avg"
$ws2.Range("A13").Formula = "=_xlfn.STDEV.S(A1:A10)"
$ws2.Range("B13").Value = "std"

$ws2.Range("B16").Select()
